$d = $word.ActiveDocument

# 1. Rename the first heading.
$d.Content.Find.Execute("Objet geoReq", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Objet geoResourceRequest", 2)

# 2. Delete the first table (the one describing "resourceRequest").
$d.Tables(1).Delete()

# 3. Delete the now-orphaned "Type resource" heading paragraph.
$d.Paragraphs(2).Range.Delete()
